# 003 Week 3 data update
# Fill in WK3 (column F) scores for week 3, plus a couple of
# WK1 (D) / WK2 (E) corrections, and clear the stray "X" placeholders
# that previously sat in column D for players who hadn't played yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 12: BAZ MASON ---
$ws.Range("F12").Value = 31

# --- Row 13: MICK SKINNER ---
$ws.Range("F13").Value = 35

# --- Row 16: LES DOBBINS ---
$ws.Range("F16").Value = 42

# --- Row 17: JOHN ANTCLIFFE ---
$ws.Range("F17").Value = 32

# --- Row 18: PAUL DIXON --- (clear "X" placeholder, add WK3 score)
$ws.Range("D18").ClearContents()
$ws.Range("F18").Value = 33

# --- Row 20: ALBIE GILLESPIE --- (clear "X" placeholder, add WK2 + WK3 scores)
$ws.Range("D20").ClearContents()
$ws.Range("E20").Value = 32
$ws.Range("F20").Value = 33

# --- Row 21: CHRIS DUFFY --- (add WK2 + WK3 scores)
$ws.Range("E21").Value = 39
$ws.Range("F21").Value = 33

# --- Row 23: TONY SLATER ---
$ws.Range("F23").Value = 30

# --- Row 24: ADY STEANE --- (clear "X" placeholder, add WK3 score)
$ws.Range("D24").ClearContents()
$ws.Range("F24").Value = 33

# --- Row 25: STEVE FELLOWS ---
$ws.Range("F25").Value = 33

# --- Row 26: MAL JONES ---
$ws.Range("F26").Value = 36

# --- Row 27: KEN PEEL ---
$ws.Range("F27").Value = 29

# --- Row 28: PAUL HANCOX --- (add WK2 + WK3 scores)
$ws.Range("E28").Value = 31
$ws.Range("F28").Value = 33

# --- Row 29: FRED HOLLIWORTH --- (add WK2 + WK3 scores)
$ws.Range("E29").Value = 30
$ws.Range("F29").Value = 32

# --- Row 30: SCOTT LEONARDE --- (clear "X" placeholder, add WK3 score)
$ws.Range("D30").ClearContents()
$ws.Range("F30").Value = 27

# --- Row 31: ANDY THOMPSON ---
$ws.Range("F31").Value = 34
